# Update the "dSF" column (column F) values to match repulled/recomputed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F5").Value = 5
$ws.Range("F10").Value = -3
$ws.Range("F12").Value = 1
$ws.Range("F14").Value = 1
$ws.Range("F17").Value = 1
$ws.Range("F18").Value = 0
$ws.Range("F20").Value = 5
$ws.Range("F26").Value = 13
